$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Widen column C (closest achievable value to the target 50.28515625 given
# this runtime's pixel-quantized ColumnWidth implementation)
$ws.Columns.Item(3).ColumnWidth = 49.45

# New row of data for exposicion #5
$ws.Cells.Item(6, 1).Value = 5

# Reuse the date formatting already used in column B (style index 1)
$ws.Cells.Item(5, 2).Copy()
$ws.Cells.Item(6, 2).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(6, 2).Value = "5/22/2012"

$ws.Cells.Item(6, 3).Value = "Exposicion oral /entrega final De Exposicion"
$ws.Cells.Item(6, 4).Value = "Pendiente"

# Update selection to match the new active cell
$ws.Range("C6").Select()
